$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44253
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 667
$ws.Range("D3").Value = 44253
$ws.Range("J3").Value = 100
$ws.Range("D4").Value = 44397
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 806
$ws.Range("D5").Value = 44342
$ws.Range("I5").Value = "Primera"
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("P5").Value = 639
$ws.Range("D6").Value = 44342
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("P6").Value = 500
$ws.Range("D7").Value = 44349
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11500
$ws.Range("P7").Value = 639
$ws.Range("D8").Value = 44349
$ws.Range("I8").Value = "Segunda"
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("P8").Value = 556
$ws.Range("D9").Value = 44356
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 11000
$ws.Range("M9").Value = 10500
$ws.Range("P9").Value = 583
$ws.Range("D10").Value = 44356
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("P10").Value = 500
$ws.Range("D11").Value = 44364
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("P11").Value = 750
$ws.Range("D12").Value = 44364
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 11000
$ws.Range("P12").Value = 611
$ws.Range("D13").Value = 44259
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("P13").Value = 694
$ws.Range("D14").Value = 44259
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("P14").Value = 556
$ws.Range("D15").Value = 44350
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 13500
$ws.Range("P15").Value = 750
$ws.Range("D16").Value = 44350
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = 11000
$ws.Range("P16").Value = 611
$ws.Range("D17").Value = 44308
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = 10500
$ws.Range("P17").Value = 583
$ws.Range("D18").Value = 44308
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 8000
$ws.Range("P18").Value = 444
$ws.Range("D19").Value = 44320
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 9000
$ws.Range("M19").Value = 9500
$ws.Range("P19").Value = 528
$ws.Range("D20").Value = 44320
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("P20").Value = 444
$ws.Range("D21").Value = 44335
$ws.Range("K21").Value = 12000
$ws.Range("M21").Value = 12500
$ws.Range("P21").Value = 694
$ws.Range("D22").Value = 44335
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 10000
$ws.Range("P22").Value = 556
$ws.Range("D23").Value = 44316
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("P23").Value = 583
$ws.Range("D24").Value = 44316
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("P24").Value = 500
$ws.Range("D25").Value = 44265
$ws.Range("I25").Value = "Primera"
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13500
$ws.Range("P25").Value = 750
$ws.Range("D26").Value = 44272
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 11000
$ws.Range("M26").Value = 10500
$ws.Range("P26").Value = 583
$ws.Range("D27").Value = 44272
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 9000
$ws.Range("P27").Value = 500
$ws.Range("D28").Value = 44313
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("P28").Value = 583
$ws.Range("D29").Value = 44313
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = 9000
$ws.Range("P29").Value = 500
$ws.Range("D30").Value = 44280
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 11000
$ws.Range("M30").Value = 10500
$ws.Range("P30").Value = 583
$ws.Range("D31").Value = 44280
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 9000
$ws.Range("P31").Value = 500
$ws.Range("D32").Value = 44328
$ws.Range("K32").Value = 9000
$ws.Range("M32").Value = 9500
$ws.Range("P32").Value = 528
$ws.Range("D33").Value = 44328
$ws.Range("J33").Value = 100
$ws.Range("D34").Value = 44379
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15500
$ws.Range("P34").Value = 861
$ws.Range("D35").Value = 44379
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 13000
$ws.Range("P35").Value = 722
$ws.Range("D36").Value = 44392
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15500
$ws.Range("P36").Value = 861
$ws.Range("D37").Value = 44392
$ws.Range("K37").Value = 14000
$ws.Range("L37").Value = 14000
$ws.Range("M37").Value = 14000
$ws.Range("P37").Value = 778
$ws.Range("D38").Value = 44384
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 16000
$ws.Range("M38").Value = 15500
$ws.Range("P38").Value = 861
$ws.Range("D39").Value = 44384
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 13000
$ws.Range("M39").Value = 13000
$ws.Range("P39").Value = 722
